$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SYSoCCtaSC")

# New fuel-type rows, added before the shared-string text for B1 is replaced so the
# shared-string table gets the same append order as the authoritative edit.
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"

# Header cell gets renamed (units clarified) and wrapped.
$ws.Range("B1").Value = "Soft cost share (dimensionless)"
$ws.Range("B1").WrapText = $true

# New column B gets an explicit width now that it holds a wrapped header.
$ws.Columns("B").ColumnWidth = 13.5

# Put the UI selection/cursor on the header cell, matching the saved view state.
$null = $ws.Range("B1").Select()

# Restore "About" as the active/visible sheet (tab) as in the source workbook.
$wsAbout = $wb.Worksheets.Item("About")
$null = $wsAbout.Activate()
$null = $wsAbout.Range("A1").Select()

Write-Host "edit applied"
